$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '257.21'
$c.Style = 'Normal'

$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '-0.46%'
$c.Style = 'Normal'

$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '27.09'
$c.Style = 'Normal'

$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '-0.04%'
$c.Style = 'Normal'

$c = $ws.Range("D4")
$c.NumberFormat = '@'
$c.Value = '4.585'
$c.Style = 'Normal'

$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '-5.67%'
$c.Style = 'Normal'

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '0.05894'
$c.Style = 'Normal'

$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '6.628'
$c.Style = 'Normal'

$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '-0.85%'
$c.Style = 'Normal'

$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.8524'
$c.Style = 'Normal'

$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '-2.61%'
$c.Style = 'Normal'

$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.9443'
$c.Style = 'Normal'

$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '-1.68%'
$c.Style = 'Normal'

$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.1392'
$c.Style = 'Normal'

$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '-1.71%'
$c.Style = 'Normal'

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.04906'
$c.Style = 'Normal'

$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '37.06%'
$c.Style = 'Normal'

$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.07083'
$c.Style = 'Normal'

$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '-2.07%'
$c.Style = 'Normal'

$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '0.03080'
$c.Style = 'Normal'

$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '-1.78%'
$c.Style = 'Normal'

$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.09129'
$c.Style = 'Normal'

$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '-1.27%'
$c.Style = 'Normal'

$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '0.001533'
$c.Style = 'Normal'

$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '-0.89%'
$c.Style = 'Normal'

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '0.0006061'
$c.Style = 'Normal'

$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '-0.12%'
$c.Style = 'Normal'

$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '0.006084'
$c.Style = 'Normal'

$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '0.98%'
$c.Style = 'Normal'

$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '3.495'
$c.Style = 'Normal'

$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '0.34%'
$c.Style = 'Normal'

$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '3.181'
$c.Style = 'Normal'

$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '-1.27%'
$c.Style = 'Normal'

$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '0.3054'
$c.Style = 'Normal'

$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '-2.86%'
$c.Style = 'Normal'

$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '-2.88%'
$c.Style = 'Normal'

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '3.955'
$c.Style = 'Normal'

$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '12.03%'
$c.Style = 'Normal'

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '0.04266'
$c.Style = 'Normal'

$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '0.48%'
$c.Style = 'Normal'

$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '0.001221'
$c.Style = 'Normal'

$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '0.05%'
$c.Style = 'Normal'

$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '0.004288'
$c.Style = 'Normal'

$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '-5.18%'
$c.Style = 'Normal'

$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '0.0001200'
$c.Style = 'Normal'

$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '0.05%'
$c.Style = 'Normal'

$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '29.79%'
$c.Style = 'Normal'

$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '0.03821'
$c.Style = 'Normal'

$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '-0.33%'
$c.Style = 'Normal'

$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '0.006131'
$c.Style = 'Normal'

$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '4.13%'
$c.Style = 'Normal'

$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.1102'
$c.Style = 'Normal'

$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '-0.15%'
$c.Style = 'Normal'

$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '0.002411'
$c.Style = 'Normal'

$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '9.60%'
$c.Style = 'Normal'

$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.01416'
$c.Style = 'Normal'

$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '31.24%'
$c.Style = 'Normal'

$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.00005360'
$c.Style = 'Normal'

$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '-2.36%'
$c.Style = 'Normal'

$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '0.00000000750'
$c.Style = 'Normal'

$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '0.05%'
$c.Style = 'Normal'

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '0.05101'
$c.Style = 'Normal'

$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '-53.23%'
$c.Style = 'Normal'

$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '11,599.75%'
$c.Style = 'Normal'

$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '0.00002100'
$c.Style = 'Normal'

$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '0.05%'
$c.Style = 'Normal'

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '0.0002000'
$c.Style = 'Normal'

$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '0.05%'
$c.Style = 'Normal'
